# New submission synced into the "JSS 3B" results sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JSS 3B")
$ws.Activate()

# Append the new row right after the last existing data row (row 10).
$newRow = 11
$ws.Cells.Item($newRow, 1).Value = "2026-02-12 21:24:31"
$ws.Cells.Item($newRow, 2).Value = "Nguru Ali "
$ws.Cells.Item($newRow, 3).Value = "Serial number 31"
$ws.Cells.Item($newRow, 4).Value = 9
